# Refresh the coin Price (column D) and Volume(1h) (column E) values
# with the latest scrape results. Column D strings get an apostrophe
# prefix so numeric-looking prices ("245.60", "41.32", ...) are stored
# as literal text (matching the sheet's existing inline-string cells)
# instead of being auto-coerced to numbers; Style is reset to "Normal"
# right after so no stray NumberFormat/quote-prefix formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.Value = "'35.365.71"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.43%  "
$dCell = $ws.Cells.Item(3, 4)
$dCell.Value = "'1.909.75"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.68%  "
$ws.Cells.Item(4, 5).Value = "  -0.56%  "
$dCell = $ws.Cells.Item(5, 4)
$dCell.Value = "'245.60"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.61%  "
$ws.Cells.Item(6, 5).Value = "  +5.94%  "
$ws.Cells.Item(7, 5).Value = "  -0.52%  "
$dCell = $ws.Cells.Item(8, 4)
$dCell.Value = "'41.32"
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.00%  "
$ws.Cells.Item(9, 5).Value = "  +6.25%  "
$dCell = $ws.Cells.Item(10, 4)
$dCell.Value = "'52.74"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +12.38%  "
$dCell = $ws.Cells.Item(11, 4)
$dCell.Value = "'0.0719"
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.73%  "
$ws.Cells.Item(12, 5).Value = "  +0.45%  "
$dCell = $ws.Cells.Item(13, 4)
$dCell.Value = "'2.187.53"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.75%  "
$ws.Cells.Item(14, 5).Value = "  +5.51%  "
$dCell = $ws.Cells.Item(15, 4)
$dCell.Value = "'0.701"
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.55%  "
$dCell = $ws.Cells.Item(16, 4)
$dCell.Value = "'1.897.31"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.98%  "
$dCell = $ws.Cells.Item(17, 4)
$dCell.Value = "'4.86"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.85%  "
$dCell = $ws.Cells.Item(18, 4)
$dCell.Value = "'35.346.17"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.37%  "
$dCell = $ws.Cells.Item(19, 4)
$dCell.Value = "'72.67"
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +3.96%  "
$dCell = $ws.Cells.Item(20, 4)
$dCell.Value = "'0.0₃0823"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.44%  "
$dCell = $ws.Cells.Item(21, 4)
$dCell.Value = "'239.62"
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.38%  "
$ws.Cells.Item(22, 5).Value = "  +2.64%  "
$dCell = $ws.Cells.Item(23, 4)
$dCell.Value = "'4.84"
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.64%  "
$ws.Cells.Item(24, 5).Value = "  -0.54%  "
$ws.Cells.Item(25, 5).Value = "  +1.15%  "
$ws.Cells.Item(26, 5).Value = "  +23.00%  "
$dCell = $ws.Cells.Item(27, 4)
$dCell.Value = "'169.85"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.39%  "
$dCell = $ws.Cells.Item(28, 4)
$dCell.Value = "'8.47"
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.74%  "
$dCell = $ws.Cells.Item(29, 4)
$dCell.Value = "'18.47"
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.61%  "
$ws.Cells.Item(30, 5).Value = "  +2.29%  "
$dCell = $ws.Cells.Item(31, 4)
$dCell.Value = "'4.15"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.76%  "
$dCell = $ws.Cells.Item(32, 4)
$dCell.Value = "'0.0566"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.29%  "
$dCell = $ws.Cells.Item(33, 4)
$dCell.Value = "'0.941"
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +15.76%  "
$ws.Cells.Item(34, 5).Value = "  -0.46%  "
$dCell = $ws.Cells.Item(35, 4)
$dCell.Value = "'4.11"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.54%  "
$ws.Cells.Item(36, 5).Value = "  -3.99%  "
$ws.Cells.Item(37, 5).Value = "  +0.47%  "
$ws.Cells.Item(38, 5).Value = "  +1.28%  "
$ws.Cells.Item(39, 5).Value = "  +1.63%  "
$ws.Cells.Item(40, 5).Value = "  +3.55%  "
$dCell = $ws.Cells.Item(41, 4)
$dCell.Value = "'16.23"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +8.30%  "
$dCell = $ws.Cells.Item(42, 4)
$dCell.Value = "'0.0642"
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +8.82%  "
$dCell = $ws.Cells.Item(43, 4)
$dCell.Value = "'90.10"
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.26%  "
$dCell = $ws.Cells.Item(44, 4)
$dCell.Value = "'1.342.30"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.29%  "
$dCell = $ws.Cells.Item(45, 4)
$dCell.Value = "'2.39"
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.73%  "
$dCell = $ws.Cells.Item(46, 4)
$dCell.Value = "'48.05"
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +38.16%  "
$ws.Cells.Item(47, 5).Value = "  +1.65%  "
$ws.Cells.Item(48, 5).Value = "  -0.88%  "
$ws.Cells.Item(49, 5).Value = "  -0.56%  "
$dCell = $ws.Cells.Item(50, 4)
$dCell.Value = "'2.094.68"
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.45%  "
$dCell = $ws.Cells.Item(51, 4)
$dCell.Value = "'0.0706"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.60%  "
